# updated GSC export data
# Append two new daily rows (2025-11-13 and 2025-11-14) to the "Chart"
# sheet's data table, and keep the "Table" sheet's header row (which
# mirrors the last-used shared strings) pointing at the same text.

$wb = $excel.ActiveWorkbook
$wsChart = $wb.Worksheets.Item(1)
$wsTable = $wb.Worksheets.Item(2)

# New rows appended at the bottom of the Chart sheet's data range.
# The date column holds plain text (not real dates), so force a text
# number format before assigning, then clear the format again so the
# cell keeps the sheet's default style (matches the other rows).
$wsChart.Range("A39").NumberFormat = "@"
$wsChart.Range("A39").Value = "2025-11-13"
$wsChart.Range("A39").ClearFormats()
$wsChart.Range("B39").Value = 0
$wsChart.Range("C39").Value = 43

$wsChart.Range("A40").NumberFormat = "@"
$wsChart.Range("A40").Value = "2025-11-14"
$wsChart.Range("A40").ClearFormats()
$wsChart.Range("B40").Value = 0
$wsChart.Range("C40").Value = 40

# The Table sheet's header row stays textually the same ("Issue",
# "Validation", "Pages") even though the underlying shared-string
# indices shift because of the two newly inserted date strings.
$wsTable.Range("A1").Value = "Issue"
$wsTable.Range("B1").Value = "Validation"
$wsTable.Range("C1").Value = "Pages"
